# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
#
# A handful of match rows in the "Iraq League" sheet were entered in the
# wrong relative order. For each affected pair of adjacent rows, the
# match id (column B) and every match-specific column (E..AD: HomeTeam,
# AwayTeam, scores, odds, etc.) need to be swapped between the two rows.
# Columns A (row index), C (Div) and D (Date) are identical for both
# rows in a pair and must stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($ws, $rowA, $rowB)

    # Swap column B (numeric match id).
    $bA = $ws.Range("B$rowA").Value()
    $bB = $ws.Range("B$rowB").Value()
    $ws.Range("B$rowA").Value = $bB
    $ws.Range("B$rowB").Value = $bA

    # Swap the contiguous block E..AD (HomeTeam through PL_AhUnder) in one
    # shot using 2D value arrays.
    $rangeA = $ws.Range("E" + $rowA + ":AD" + $rowA)
    $rangeB = $ws.Range("E" + $rowB + ":AD" + $rowB)

    $valA = $rangeA.Value()
    $valB = $rangeB.Value()

    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

$pairs = @(
    @(17, 18),
    @(73, 74),
    @(78, 79),
    @(108, 109),
    @(135, 136),
    @(219, 220),
    @(223, 224)
)

foreach ($pair in $pairs) {
    Swap-Rows $ws $pair[0] $pair[1]
}
